# Add two new columns (I, J) with headers "I0" and "IF" to match the
# existing header row formatting, plus the corresponding data values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells - copy the style used by the other header cells (e.g. H1)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats

# New data cells for rows 2 and 3
$ws.Range("I2").Value = 8
$ws.Range("J2").Value = 8
$ws.Range("I3").Value = 8
$ws.Range("J3").Value = 8
